$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2431.942902
$ws.Range("K2").Value = 2421.7722246
$ws.Range("L2").Value = 2411.6015472
$ws.Range("M2").Value = 2401.4308698
$ws.Range("N2").Value = 2391.2601924
$ws.Range("O2").Value = 2381.089515
$ws.Range("P2").Value = 2420.084159
$ws.Range("Q2").Value = 2459.078803
$ws.Range("R2").Value = 2498.073447
$ws.Range("S2").Value = 2537.068091
$ws.Range("T2").Value = 2576.062735
$ws.Range("U2").Value = 2644.5606574
$ws.Range("V2").Value = 2713.0585798
$ws.Range("W2").Value = 2781.5565022
$ws.Range("X2").Value = 2850.0544246
$ws.Range("Y2").Value = 2918.552347
$ws.Range("Z2").Value = 2984.2074356
$ws.Range("AA2").Value = 3049.8625242
$ws.Range("AB2").Value = 3115.5176128
$ws.Range("AC2").Value = 3181.1727014
$ws.Range("AD2").Value = 3246.82779
$ws.Range("AE2").Value = 3314.4815632
$ws.Range("AF2").Value = 3382.1353364
$ws.Range("AG2").Value = 3449.7891096
$ws.Range("AH2").Value = 3517.4428828
$ws.Range("AI2").Value = 3585.096656
$ws.Range("AJ2").Value = 3653.2644508
$ws.Range("AK2").Value = 3721.4322456
$ws.Range("AL2").Value = 3789.6000404
$ws.Range("AM2").Value = 3857.7678352
$ws.Range("AN2").Value = 3925.93563
$ws.Range("AO2").Value = 3993.1127378
$ws.Range("AP2").Value = 4060.2898456
$ws.Range("AQ2").Value = 4127.4669534
$ws.Range("AR2").Value = 4194.6440612
$ws.Range("AS2").Value = 4261.821169

$ws.Range("J8:AS8").Value = -0.1

$ws.Range("J9:AS9").Value = 1

$ws.Range("J13:AS13").Value = 3.145207224
